$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 732
$ws.Range("F3").Value = 14024
$ws.Range("F4").Value = 14024
$ws.Range("F5").Value = 14033
$ws.Range("F6").Value = 1359
$ws.Range("F7").Value = 1380
$ws.Range("F8").Value = 5808
$ws.Range("F9").Value = 970
$ws.Range("F10").Value = 566
$ws.Range("F11").Value = 113
$ws.Range("F12").Value = 376
$ws.Range("F14").Value = 1521
$ws.Range("F15").Value = 421
$ws.Range("F16").Value = 2131
$ws.Range("F17").Value = 1175
$ws.Range("F18").Value = 1784
$ws.Range("F19").Value = 908
$ws.Range("F20").Value = 30
$ws.Range("F21").Value = 2249
$ws.Range("F22").Value = 549
$ws.Range("F23").Value = 789
$ws.Range("F24").Value = 3264
$ws.Range("F26").Value = 304
$ws.Range("F27").Value = 2321
$ws.Range("F28").Value = 74
$ws.Range("F29").Value = 115
$ws.Range("F31").Value = 1761
$ws.Range("F32").Value = 1063
$ws.Range("F33").Value = 1339
$ws.Range("F34").Value = 91
$ws.Range("F35").Value = 134
$ws.Range("F36").Value = 4653
$ws.Range("F37").Value = 4733
$ws.Range("F38").Value = 292
$ws.Range("F39").Value = 155
$ws.Range("F40").Value = 662
$ws.Range("F41").Value = 672
$ws.Range("F42").Value = 3258
$ws.Range("F43").Value = 41
$ws.Range("F45").Value = 328
$ws.Range("F46").Value = 83
$ws.Range("F48").Value = 4404
$ws.Range("F49").Value = 533
$ws.Range("F50").Value = 275

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 110
$ws.Range("F7").Value = 80
$ws.Range("F17").Value = 16
$ws.Range("F22").Value = 53
$ws.Range("F24").Value = 1

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 7408
$ws.Range("F3").Value = 213
$ws.Range("F4").Value = 669

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 7408
$ws.Range("F3").Value = 732
$ws.Range("F4").Value = 213
$ws.Range("F5").Value = 669
$ws.Range("F7").Value = 14024
$ws.Range("F8").Value = 14024
$ws.Range("F9").Value = 14033
$ws.Range("F10").Value = 1359
$ws.Range("F11").Value = 1380
$ws.Range("F12").Value = 5808
$ws.Range("F13").Value = 970
$ws.Range("F15").Value = 80
$ws.Range("F16").Value = 421
$ws.Range("F17").Value = 1175
$ws.Range("F18").Value = 1785
$ws.Range("F20").Value = 789
$ws.Range("F21").Value = 3264
$ws.Range("F22").Value = 304
$ws.Range("F23").Value = 74
$ws.Range("F24").Value = 115
$ws.Range("F26").Value = 1761
$ws.Range("F31").Value = 16
$ws.Range("F32").Value = 1063
$ws.Range("F33").Value = 1339
$ws.Range("F34").Value = 91
$ws.Range("F36").Value = 4654
$ws.Range("F37").Value = 4733
$ws.Range("F38").Value = 292
$ws.Range("F39").Value = 155
$ws.Range("F40").Value = 3258
$ws.Range("F41").Value = 41
$ws.Range("F43").Value = 328
$ws.Range("F44").Value = 83
$ws.Range("F46").Value = 4404
$ws.Range("F47").Value = 275
